$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 22.11125233333334
$ws.Range("H2").Value = 66.33375700000001
$ws.Range("I2").Value = 0.08763778737242772
$ws.Range("J2").Value = 0.08763778737242772
$ws.Range("M2").Value = 14.440165
$ws.Range("N2").Value = 43.320495
$ws.Range("O2").Value = 0.1441015470002482
$ws.Range("P2").Value = 0.1441015470002482
$ws.Range("Q2").Value = 319.2901320499684
$ws.Range("R2").Value = 2873.611188449715
$ws.Range("S2").Value = 0.01262874073604566
$ws.Range("T2").Value = 0.01262874073604565

# Row 3
$ws.Range("G3").Value = 22.11125233333334
$ws.Range("H3").Value = 66.33375700000001
$ws.Range("I3").Value = 0.08763778737242772
$ws.Range("J3").Value = 0.08763778737242772
$ws.Range("O3").Value = 0.3846359116098663
$ws.Range("P3").Value = 0.3846359116098662
$ws.Range("Q3").Value = 852.249358633621
$ws.Range("R3").Value = 7670.244227702588
$ws.Range("S3").Value = 0.03370864023746536
$ws.Range("T3").Value = 0.03370864023746536

# Row 4
$ws.Range("G4").Value = 22.11125233333334
$ws.Range("H4").Value = 66.33375700000001
$ws.Range("I4").Value = 0.08763778737242772
$ws.Range("J4").Value = 0.08763778737242772
$ws.Range("M4").Value = 21.954262
$ws.Range("N4").Value = 65.862786
$ws.Range("O4").Value = 0.2190863551385157
$ws.Range("P4").Value = 0.2190863551385156
$ws.Range("Q4").Value = 485.4362268741114
$ws.Range("R4").Value = 4368.926041867002
$ws.Range("S4").Value = 0.01920024340782942
$ws.Range("T4").Value = 0.01920024340782942

# Row 5
$ws.Range("G5").Value = 22.11125233333334
$ws.Range("H5").Value = 66.33375700000001
$ws.Range("I5").Value = 0.08763778737242772
$ws.Range("J5").Value = 0.08763778737242772
$ws.Range("M5").Value = 25.27013633333333
$ws.Range("N5").Value = 75.81040899999999
$ws.Range("O5").Value = 0.2521761862513699
$ws.Range("P5").Value = 0.2521761862513699
$ws.Range("Q5").Value = 558.7543609640682
$ws.Range("R5").Value = 5028.789248676613
$ws.Range("S5").Value = 0.02210016299108729
$ws.Range("T5").Value = 0.02210016299108729

# Row 6
$ws.Range("I6").Value = 0.8254813202458152
$ws.Range("J6").Value = 0.8254813202458152
$ws.Range("M6").Value = 14.440165
$ws.Range("N6").Value = 43.320495
$ws.Range("O6").Value = 0.1441015470002482
$ws.Range("P6").Value = 0.1441015470002482
$ws.Range("Q6").Value = 3007.470266518748
$ws.Range("R6").Value = 27067.23239866873
$ws.Range("S6").Value = 0.1189531352672293
$ws.Range("T6").Value = 0.1189531352672293

# Row 7
$ws.Range("I7").Value = 0.8254813202458152
$ws.Range("J7").Value = 0.8254813202458152
$ws.Range("O7").Value = 0.3846359116098663
$ws.Range("P7").Value = 0.3846359116098662
$ws.Range("S7").Value = 0.3175097601296651
$ws.Range("T7").Value = 0.317509760129665

# Row 8
$ws.Range("I8").Value = 0.8254813202458152
$ws.Range("J8").Value = 0.8254813202458152
$ws.Range("M8").Value = 21.954262
$ws.Range("N8").Value = 65.862786
$ws.Range("O8").Value = 0.2190863551385157
$ws.Range("P8").Value = 0.2190863551385156
$ws.Range("Q8").Value = 4572.440147904294
$ws.Range("R8").Value = 41151.96133113866
$ws.Range("S8").Value = 0.1808516936875854
$ws.Range("T8").Value = 0.1808516936875854

# Row 9
$ws.Range("I9").Value = 0.8254813202458152
$ws.Range("J9").Value = 0.8254813202458152
$ws.Range("M9").Value = 25.27013633333333
$ws.Range("N9").Value = 75.81040899999999
$ws.Range("O9").Value = 0.2521761862513699
$ws.Range("P9").Value = 0.2521761862513699
$ws.Range("Q9").Value = 5263.041222408129
$ws.Range("R9").Value = 47367.37100167317
$ws.Range("S9").Value = 0.2081667311613355
$ws.Range("T9").Value = 0.2081667311613355

# Row 10
$ws.Range("G10").Value = 12.43397833333333
$ws.Range("H10").Value = 37.301935
$ws.Range("I10").Value = 0.04928198244688778
$ws.Range("J10").Value = 0.04928198244688778
$ws.Range("M10").Value = 14.440165
$ws.Range("N10").Value = 43.320495
$ws.Range("O10").Value = 0.1441015470002482
$ws.Range("P10").Value = 0.1441015470002482
$ws.Range("Q10").Value = 179.5486987397584
$ws.Range("R10").Value = 1615.938288657825
$ws.Range("S10").Value = 0.007101609909835608
$ws.Range("T10").Value = 0.007101609909835607

# Row 11
$ws.Range("G11").Value = 12.43397833333333
$ws.Range("H11").Value = 37.301935
$ws.Range("I11").Value = 0.04928198244688778
$ws.Range("J11").Value = 0.04928198244688778
$ws.Range("O11").Value = 0.3846359116098663
$ws.Range("P11").Value = 0.3846359116098662
$ws.Range("Q11").Value = 479.2514643719489
$ws.Range("R11").Value = 4313.26317934754
$ws.Range("S11").Value = 0.01895562024440011
$ws.Range("T11").Value = 0.01895562024440011

# Row 12
$ws.Range("G12").Value = 12.43397833333333
$ws.Range("H12").Value = 37.301935
$ws.Range("I12").Value = 0.04928198244688778
$ws.Range("J12").Value = 0.04928198244688778
$ws.Range("M12").Value = 21.954262
$ws.Range("N12").Value = 65.862786
$ws.Range("O12").Value = 0.2190863551385157
$ws.Range("P12").Value = 0.2190863551385156
$ws.Range("Q12").Value = 272.9788180323234
$ws.Range("R12").Value = 2456.80936229091
$ws.Range("S12").Value = 0.01079700990828895
$ws.Range("T12").Value = 0.01079700990828895

# Row 13
$ws.Range("G13").Value = 12.43397833333333
$ws.Range("H13").Value = 37.301935
$ws.Range("I13").Value = 0.04928198244688778
$ws.Range("J13").Value = 0.04928198244688778
$ws.Range("M13").Value = 25.27013633333333
$ws.Range("N13").Value = 75.81040899999999
$ws.Range("O13").Value = 0.2521761862513699
$ws.Range("P13").Value = 0.2521761862513699
$ws.Range("Q13").Value = 314.2083276490461
$ws.Range("R13").Value = 2827.874948841415
$ws.Range("S13").Value = 0.01242774238436312
$ws.Range("T13").Value = 0.01242774238436312

# Row 14
$ws.Range("G14").Value = 9.486307333333334
$ws.Range("H14").Value = 28.458922
$ws.Range("I14").Value = 0.03759890993486929
$ws.Range("J14").Value = 0.03759890993486929
$ws.Range("M14").Value = 14.440165
$ws.Range("N14").Value = 43.320495
$ws.Range("O14").Value = 0.1441015470002482
$ws.Range("P14").Value = 0.1441015470002482
$ws.Range("Q14").Value = 136.9838431340434
$ws.Range("R14").Value = 1232.85458820639
$ws.Range("S14").Value = 0.005418061087137668
$ws.Range("T14").Value = 0.005418061087137667

# Row 15
$ws.Range("G15").Value = 9.486307333333334
$ws.Range("H15").Value = 28.458922
$ws.Range("I15").Value = 0.03759890993486929
$ws.Range("J15").Value = 0.03759890993486929
$ws.Range("O15").Value = 0.3846359116098663
$ws.Range("P15").Value = 0.3846359116098662
$ws.Range("Q15").Value = 365.6373333701609
$ws.Range("R15").Value = 3290.736000331448
$ws.Range("S15").Value = 0.01446189099833571
$ws.Range("T15").Value = 0.01446189099833571

# Row 16
$ws.Range("G16").Value = 9.486307333333334
$ws.Range("H16").Value = 28.458922
$ws.Range("I16").Value = 0.03759890993486929
$ws.Range("J16").Value = 0.03759890993486929
$ws.Range("M16").Value = 21.954262
$ws.Range("N16").Value = 65.862786
$ws.Range("O16").Value = 0.2190863551385157
$ws.Range("P16").Value = 0.2190863551385156
$ws.Range("Q16").Value = 208.2648766085213
$ws.Range("R16").Value = 1874.383889476692
$ws.Range("S16").Value = 0.008237408134811839
$ws.Range("T16").Value = 0.008237408134811837

# Row 17
$ws.Range("G17").Value = 9.486307333333334
$ws.Range("H17").Value = 28.458922
$ws.Range("I17").Value = 0.03759890993486929
$ws.Range("J17").Value = 0.03759890993486929
$ws.Range("M17").Value = 25.27013633333333
$ws.Range("N17").Value = 75.81040899999999
$ws.Range("O17").Value = 0.2521761862513699
$ws.Range("P17").Value = 0.2521761862513699
$ws.Range("Q17").Value = 239.7202796132331
$ws.Range("R17").Value = 2157.482516519098
$ws.Range("S17").Value = 0.009481549714584083
$ws.Range("T17").Value = 0.009481549714584083
